$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feedback entry row (row 3)
$ws.Range("A3").Value = "V1.17"
$ws.Range("B3").Value = "- Aesthetic UI improvements`n- Models for each gamemode (don't reuse chess models for checkers)"
$ws.Range("C3").Value = "- These are aesthetic suggestions so they are a low priority and will be sorted later"

# B3/C3 pick up the same wrapped/quote-prefixed format already used by B2/C2
$ws.Range("B2:C2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)

# Whole-sheet wrap text (matches the header/version columns getting the new wrap style)
$ws.Columns.Item(1).WrapText = $true
$ws.Columns.Item(2).WrapText = $true
$ws.Columns.Item(3).WrapText = $true

$ws.Rows.Item(3).RowHeight = 90

$ws.Range("D3").Select()
